# CS 513 Homework 2 - add a new "Slope Evaluation" slide at the end of the deck.
$p = $ppt.ActivePresentation

# Append a new slide after the existing ones, using the same "Title and
# Content" layout already used by slides 2 and 3 (ppLayoutText = 2, which
# PowerPoint maps onto the deck's "Title and Content" custom layout).
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Slope Evaluation"

# Body / content placeholder.
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Our slope calculations "
